{"js": "// Remove the \"Date\" styled paragraph (e.g. \"2023-01-27\") that follows the\n// title paragraph \"Helloworld\", per the commit removing the R Markdown\n// date line.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"style,text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.style === \"Date\") {\n    p.delete();\n  }\n}\nawait context.sync();\n", "ps1": "# Remove the \"Date\" styled paragraph (e.g. \"2023-01-27\") that follows the\n# title paragraph \"Helloworld\", per the commit removing the R Markdown\n# date line.\n$d = $word.ActiveDocument\n\n# Walk backwards so deleting a paragraph doesn't perturb indices of the\n# ones we still need to visit.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Style.NameLocal -eq \"Date\") {\n        $p.Range.Delete()\n    }\n}\n"}
